# Auto-generated Excel COM-interop script to apply the market-data refresh diff.
# Updates cached profit-calculation columns (H-N) across ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets
# to reflect a scheduled market-data resync. Values come from the authoritative diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2499
$ws.Cells.Item(40, 9).Value = 2499
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 2499
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = -2324
$ws.Cells.Item(40, 14).ClearContents()

$ws.Cells.Item(45, 8).Value = 200
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 10).Value = 200
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 12).Value = 600
$ws.Cells.Item(45, 14).Value = -984

$ws.Cells.Item(49, 9).Value = 0
$ws.Cells.Item(49, 10).Value = 500
$ws.Cells.Item(49, 11).Value = 0
$ws.Cells.Item(49, 12).Value = 1500
$ws.Cells.Item(49, 13).Value = -1772

$ws.Cells.Item(51, 8).Value = 36665
$ws.Cells.Item(51, 9).Value = 10000
$ws.Cells.Item(51, 10).Value = 49997.5
$ws.Cells.Item(51, 11).Value = 10000
$ws.Cells.Item(51, 12).Value = 49997.5
$ws.Cells.Item(51, 13).Value = -9516
$ws.Cells.Item(51, 14).Value = -50965.5

$ws.Cells.Item(62, 8).Value = 5928.5625
$ws.Cells.Item(62, 9).Value = 2981.6
$ws.Cells.Item(62, 10).Value = 7268.091
$ws.Cells.Item(62, 11).Value = 2981.6
$ws.Cells.Item(62, 12).Value = 7268.091
$ws.Cells.Item(62, 13).Value = -2357.6
$ws.Cells.Item(62, 14).Value = -8516.091

$ws.Cells.Item(65, 8).Value = 5928.5625
$ws.Cells.Item(65, 9).Value = 2981.6
$ws.Cells.Item(65, 10).Value = 7268.091
$ws.Cells.Item(65, 11).Value = 14908
$ws.Cells.Item(65, 12).Value = 36340.455
$ws.Cells.Item(65, 13).Value = -11788
$ws.Cells.Item(65, 14).Value = -42580.455

$ws.Cells.Item(75, 8).Value = 49375
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 49375
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).Value = 49375
$ws.Cells.Item(75, 14).Value = -51247

$ws.Cells.Item(78, 8).Value = 49375
$ws.Cells.Item(78, 9).Value = 0
$ws.Cells.Item(78, 10).Value = 49375
$ws.Cells.Item(78, 11).Value = 0
$ws.Cells.Item(78, 12).Value = 148125
$ws.Cells.Item(78, 14).Value = -157485

$ws.Cells.Item(93, 8).Value = 60000
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 60000
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = 60000
$ws.Cells.Item(93, 14).Value = -64992

$ws.Cells.Item(97, 8).Value = 2361
$ws.Cells.Item(97, 9).Value = 8420
$ws.Cells.Item(97, 10).Value = 1603.625
$ws.Cells.Item(97, 11).Value = 25260
$ws.Cells.Item(97, 12).Value = 4810.875
$ws.Cells.Item(97, 13).Value = -24764
$ws.Cells.Item(97, 14).Value = -5802.875

$ws.Cells.Item(107, 8).Value = 2350.5
$ws.Cells.Item(107, 9).Value = 2350.5
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 2350.5
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = -430.5

$ws.Cells.Item(112, 8).Value = 3680.5334
$ws.Cells.Item(112, 9).Value = 4200
$ws.Cells.Item(112, 10).Value = 3550.6667
$ws.Cells.Item(112, 11).Value = 12600
$ws.Cells.Item(112, 12).Value = 10652.0001
$ws.Cells.Item(112, 13).Value = -11492
$ws.Cells.Item(112, 14).Value = -12868.0001

$ws.Cells.Item(116, 8).Value = 7000
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 10).Value = 7000
$ws.Cells.Item(116, 11).Value = 0
$ws.Cells.Item(116, 12).Value = 7000
$ws.Cells.Item(116, 14).Value = -13884
$ws.Cells.Item(116, 13).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(10, 8).Value = 6756.6665
$ws.Cells.Item(10, 9).Value = 20000
$ws.Cells.Item(10, 10).Value = 135
$ws.Cells.Item(10, 11).Value = 20000
$ws.Cells.Item(10, 12).Value = 135
$ws.Cells.Item(10, 13).Value = -19830
$ws.Cells.Item(10, 14).Value = -475

$ws.Cells.Item(23, 8).Value = 500
$ws.Cells.Item(23, 9).Value = 500
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 11).Value = 500
$ws.Cells.Item(23, 12).Value = 0
$ws.Cells.Item(23, 13).Value = -241

$ws.Cells.Item(97, 8).Value = 654.1429000000001
$ws.Cells.Item(97, 9).Value = 471.5
$ws.Cells.Item(97, 10).Value = 1750
$ws.Cells.Item(97, 11).Value = 471.5
$ws.Cells.Item(97, 12).Value = 1750
$ws.Cells.Item(97, 13).Value = 24.5
$ws.Cells.Item(97, 14).Value = -2742

$ws.Cells.Item(122, 8).Value = 4800
$ws.Cells.Item(122, 9).Value = 4000
$ws.Cells.Item(122, 10).Value = 5000
$ws.Cells.Item(122, 11).Value = 12000
$ws.Cells.Item(122, 12).Value = 15000
$ws.Cells.Item(122, 13).Value = -9550
$ws.Cells.Item(122, 14).Value = -19900

$ws.Cells.Item(132, 8).Value = 2143.0715
$ws.Cells.Item(132, 9).Value = 1533.8096
$ws.Cells.Item(132, 10).Value = 3970.8572
$ws.Cells.Item(132, 11).Value = 4601.4288
$ws.Cells.Item(132, 12).Value = 11912.5716
$ws.Cells.Item(132, 13).Value = -2071.4288
$ws.Cells.Item(132, 14).Value = -16972.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 836.3333
$ws.Cells.Item(94, 9).Value = 754.5
$ws.Cells.Item(94, 10).Value = 1000
$ws.Cells.Item(94, 11).Value = 754.5
$ws.Cells.Item(94, 12).Value = 1000
$ws.Cells.Item(94, 13).Value = -303.5
$ws.Cells.Item(94, 14).Value = -1902

$ws.Cells.Item(99, 8).Value = 1420.2106
$ws.Cells.Item(99, 9).Value = 961.5
$ws.Cells.Item(99, 10).Value = 3866.6667
$ws.Cells.Item(99, 11).Value = 961.5
$ws.Cells.Item(99, 12).Value = 3866.6667
$ws.Cells.Item(99, 13).Value = 536.5
$ws.Cells.Item(99, 14).Value = -6862.6667

$ws.Cells.Item(134, 8).Value = 2504.7144
$ws.Cells.Item(134, 9).Value = 2234.6667
$ws.Cells.Item(134, 10).Value = 4125
$ws.Cells.Item(134, 11).Value = 6704.000100000001
$ws.Cells.Item(134, 12).Value = 12375
$ws.Cells.Item(134, 13).Value = -4169.000100000001
$ws.Cells.Item(134, 14).Value = -17445

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 115401.555
$ws.Cells.Item(22, 9).Value = 170481.5
$ws.Cells.Item(22, 10).Value = 5241.6665
$ws.Cells.Item(22, 11).Value = 170481.5
$ws.Cells.Item(22, 12).Value = 5241.6665
$ws.Cells.Item(22, 13).Value = -170131.5
$ws.Cells.Item(22, 14).Value = -5941.6665

$ws.Cells.Item(122, 8).Value = 2000
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 10).Value = 2000
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 12).Value = 6000
$ws.Cells.Item(122, 14).Value = -10900
$ws.Cells.Item(122, 13).ClearContents()

$ws.Cells.Item(132, 8).Value = 4096.2607
$ws.Cells.Item(132, 9).Value = 2315.4285
$ws.Cells.Item(132, 10).Value = 6866.4443
$ws.Cells.Item(132, 11).Value = 6946.2855
$ws.Cells.Item(132, 12).Value = 20599.3329
$ws.Cells.Item(132, 13).Value = -4416.2855
$ws.Cells.Item(132, 14).Value = -25659.3329

$ws.Cells.Item(134, 8).Value = 1706.3334
$ws.Cells.Item(134, 9).Value = 1322.7
$ws.Cells.Item(134, 10).Value = 3624.5
$ws.Cells.Item(134, 11).Value = 3968.1
$ws.Cells.Item(134, 12).Value = 10873.5
$ws.Cells.Item(134, 13).Value = -1433.1
$ws.Cells.Item(134, 14).Value = -15943.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 48128.523
$ws.Cells.Item(33, 9).Value = 591.3889
$ws.Cells.Item(33, 10).Value = 333351.34
$ws.Cells.Item(33, 11).Value = 3548.3334
$ws.Cells.Item(33, 12).Value = 2000108.04
$ws.Cells.Item(33, 13).Value = -3265.3334
$ws.Cells.Item(33, 14).Value = -2000674.04

$ws.Cells.Item(35, 8).Value = 650
$ws.Cells.Item(35, 9).Value = 500
$ws.Cells.Item(35, 10).Value = 800
$ws.Cells.Item(35, 11).Value = 1500
$ws.Cells.Item(35, 12).Value = 2400
$ws.Cells.Item(35, 13).Value = -1212
$ws.Cells.Item(35, 14).Value = -2976

$ws.Cells.Item(92, 8).Value = 698.5
$ws.Cells.Item(92, 9).Value = 598
$ws.Cells.Item(92, 10).Value = 799
$ws.Cells.Item(92, 11).Value = 1794
$ws.Cells.Item(92, 12).Value = 2397
$ws.Cells.Item(92, 13).Value = -546
$ws.Cells.Item(92, 14).Value = -4893

$ws.Cells.Item(102, 8).Value = 7000
$ws.Cells.Item(102, 9).Value = 7000
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 21000
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = -18566
$ws.Cells.Item(102, 14).ClearContents()

$ws.Cells.Item(106, 8).Value = 5500
$ws.Cells.Item(106, 9).Value = 0
$ws.Cells.Item(106, 10).Value = 5500
$ws.Cells.Item(106, 11).Value = 0
$ws.Cells.Item(106, 12).Value = 16500
$ws.Cells.Item(106, 14).Value = -18392

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(27, 8).Value = 18500
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 10).Value = 18500
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 18500
$ws.Cells.Item(27, 14).Value = -18832

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(34, 8).Value = 38749.25
$ws.Cells.Item(34, 9).Value = 42499
$ws.Cells.Item(34, 10).Value = 34999.5
$ws.Cells.Item(34, 11).Value = 42499
$ws.Cells.Item(34, 12).Value = 34999.5
$ws.Cells.Item(34, 13).Value = -42327
$ws.Cells.Item(34, 14).Value = -35343.5

$ws.Cells.Item(61, 8).Value = 2249.4285
$ws.Cells.Item(61, 9).Value = 2249.4285
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 2249.4285
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -2047.4285
$ws.Cells.Item(61, 14).ClearContents()

$ws.Cells.Item(68, 8).Value = 3789.5
$ws.Cells.Item(68, 9).Value = 3499.4119
$ws.Cells.Item(68, 10).Value = 5433.3335
$ws.Cells.Item(68, 11).Value = 3499.4119
$ws.Cells.Item(68, 12).Value = 5433.3335
$ws.Cells.Item(68, 13).Value = -2750.4119
$ws.Cells.Item(68, 14).Value = -6931.3335

$ws.Cells.Item(71, 8).Value = 3789.5
$ws.Cells.Item(71, 9).Value = 3499.4119
$ws.Cells.Item(71, 10).Value = 5433.3335
$ws.Cells.Item(71, 11).Value = 17497.0595
$ws.Cells.Item(71, 12).Value = 27166.6675
$ws.Cells.Item(71, 13).Value = -13753.0595
$ws.Cells.Item(71, 14).Value = -34654.6675

$ws.Cells.Item(113, 8).Value = 2249.4285
$ws.Cells.Item(113, 9).Value = 2249.4285
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 2249.4285
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = -79.42849999999999
$ws.Cells.Item(113, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 347001.34
$ws.Cells.Item(132, 9).Value = 500502
$ws.Cells.Item(132, 10).Value = 40000
$ws.Cells.Item(132, 11).Value = 1501506
$ws.Cells.Item(132, 12).Value = 120000
$ws.Cells.Item(132, 13).Value = -1498976
$ws.Cells.Item(132, 14).Value = -125060

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 3751.5
$ws.Cells.Item(96, 9).Value = 3751.5
$ws.Cells.Item(96, 10).Value = 0
$ws.Cells.Item(96, 11).Value = 3751.5
$ws.Cells.Item(96, 12).Value = 0
$ws.Cells.Item(96, 13).Value = -2378.5

$ws.Cells.Item(122, 8).Value = 1351.1875
$ws.Cells.Item(122, 9).Value = 1351.1875
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 4053.5625
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -1603.5625

$ws.Cells.Item(132, 8).Value = 1713.579
$ws.Cells.Item(132, 9).Value = 1656.3529
$ws.Cells.Item(132, 10).Value = 2200
$ws.Cells.Item(132, 11).Value = 4969.0587
$ws.Cells.Item(132, 12).Value = 6600
$ws.Cells.Item(132, 13).Value = -2439.0587
$ws.Cells.Item(132, 14).Value = -11660

$ws.Cells.Item(136, 8).Value = 13794
$ws.Cells.Item(136, 9).Value = 13794
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 41382
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -38832
